$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.860.51'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.934.79'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '374.47'
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.44'
$ws.Range("E6").Value = '  -3.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.534'
$ws.Range("E7").Value = '  -1.12%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").Value = '  -2.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.13'
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.404.81'
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.84'
$ws.Range("E14").Value = '  -3.36%  '
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.927.82'
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.977'
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.27'
$ws.Range("E18").Value = '  +39.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '50.843.00'
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.13'
$ws.Range("E20").Value = '  -5.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.57'
$ws.Range("E21").Value = '  -2.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0951'
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '263.37'
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.34'
$ws.Range("E24").Value = '  -1.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.03'
$ws.Range("E25").Value = '  +7.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.23'
$ws.Range("E26").Value = '  +7.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.87'
$ws.Range("E27").Value = '  +6.97%  '
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.111'
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.53'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.85'
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.64'
$ws.Range("E33").Value = '  -1.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.45'
$ws.Range("E34").Value = '  -3.79%  '
$ws.Range("E35").Value = '  -3.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0441'
$ws.Range("E36").Value = '  -1.44%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.02'
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.115'
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.52'
$ws.Range("E40").Value = '  -2.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.29'
$ws.Range("E41").Value = '  -5.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.78'
$ws.Range("E42").Value = '  -3.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.24'
$ws.Range("E43").Value = '  -2.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.288'
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.00'
$ws.Range("E45").Value = '  -5.05%  '
$ws.Range("E46").Value = '  -2.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.27'
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("E48").Value = '  -3.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.986.77'
$ws.Range("E49").Value = '  -2.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0339'
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("E51").Value = '  -1.88%  '
